# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each class sheet to
# reflect the latest Universalis snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value2 = 7850
$ws.Range("I29").Value2 = 7250
$ws.Range("K29").Value2 = 21750
$ws.Range("M29").Value2 = -21469
$ws.Range("H31").Value2 = 107.5
$ws.Range("I31").Value2 = 107.5
$ws.Range("K31").Value2 = 322.5
$ws.Range("M31").Value2 = -92.5
$ws.Range("H32").Value2 = 6026.143
$ws.Range("J32").Value2 = 6738.8335
$ws.Range("L32").Value2 = 6738.8335
$ws.Range("N32").Value2 = -7390.8335
$ws.Range("H43").Value2 = 3912.5
$ws.Range("I43").Value2 = 3250
$ws.Range("K43").Value2 = 3250
$ws.Range("M43").Value2 = -3181
$ws.Range("H51").Value2 = 3131.5557
$ws.Range("I51").Value2 = 1502.3684
$ws.Range("K51").Value2 = 1502.3684
$ws.Range("M51").Value2 = -1018.3684
$ws.Range("H52").Value2 = 0
$ws.Range("I52").Value2 = 0
$ws.Range("K52").Value2 = 0
$ws.Range("M52").ClearContents()
$ws.Range("H55").Value2 = 213.64706
$ws.Range("I55").Value2 = 146.6
$ws.Range("K55").Value2 = 146.6
$ws.Range("M55").Value2 = 67.40000000000001
$ws.Range("H76").Value2 = 8205.308000000001
$ws.Range("I76").Value2 = 7037.2
$ws.Range("J76").Value2 = 12099
$ws.Range("K76").Value2 = 7037.2
$ws.Range("L76").Value2 = 12099
$ws.Range("M76").Value2 = -6722.2
$ws.Range("N76").Value2 = -12729
$ws.Range("H79").Value2 = 8205.308000000001
$ws.Range("I79").Value2 = 7037.2
$ws.Range("J79").Value2 = 12099
$ws.Range("K79").Value2 = 7037.2
$ws.Range("L79").Value2 = 12099
$ws.Range("M79").Value2 = -5945.2
$ws.Range("N79").Value2 = -14283
$ws.Range("H127").Value2 = 2247
$ws.Range("I127").Value2 = 987.25
$ws.Range("J127").Value2 = 3086.8333
$ws.Range("K127").Value2 = 2961.75
$ws.Range("L127").Value2 = 9260.499899999999
$ws.Range("M127").Value2 = 1998.25
$ws.Range("N127").Value2 = -19180.4999
$ws.Range("H129").Value2 = 17808.076
$ws.Range("I129").Value2 = 960
$ws.Range("K129").Value2 = 2880
$ws.Range("M129").Value2 = 2120
$ws.Range("H135").Value2 = 11906401
$ws.Range("I135").Value2 = 1378.4166
$ws.Range("K135").Value2 = 12405.7494
$ws.Range("M135").Value2 = -9870.749400000001
$ws.Range("H137").Value2 = 33335226
$ws.Range("I137").Value2 = 19232690
$ws.Range("K137").Value2 = 57698070
$ws.Range("M137").Value2 = -57695520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value2 = 20000000
$ws.Range("I11").Value2 = 20000000
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 20000000
$ws.Range("L11").Value2 = 0
$ws.Range("M11").Value2 = -19999856
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value2 = 16963228
$ws.Range("I32").Value2 = 29418368
$ws.Range("K32").Value2 = 29418368
$ws.Range("M32").Value2 = -29418081
$ws.Range("H74").Value2 = 45507820
$ws.Range("I74").Value2 = 47674384
$ws.Range("K74").Value2 = 47674384
$ws.Range("M74").Value2 = -47673510
$ws.Range("H77").Value2 = 45507820
$ws.Range("I77").Value2 = 47674384
$ws.Range("K77").Value2 = 238371920
$ws.Range("M77").Value2 = -238367552
$ws.Range("H80").Value2 = 75108.5
$ws.Range("I80").Value2 = 75108.5
$ws.Range("K80").Value2 = 75108.5
$ws.Range("M80").Value2 = -74110.5
$ws.Range("H83").Value2 = 75108.5
$ws.Range("I83").Value2 = 75108.5
$ws.Range("K83").Value2 = 225325.5
$ws.Range("M83").Value2 = -220333.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 341.86206
$ws.Range("I7").Value2 = 108.92308
$ws.Range("J7").Value2 = 531.125
$ws.Range("K7").Value2 = 108.92308
$ws.Range("L7").Value2 = 531.125
$ws.Range("M7").Value2 = 4.076920000000001
$ws.Range("N7").Value2 = -757.125
$ws.Range("H22").Value2 = 8722.416999999999
$ws.Range("J22").Value2 = 795.25
$ws.Range("L22").Value2 = 795.25
$ws.Range("N22").Value2 = -1495.25
$ws.Range("H31").Value2 = 21743362
$ws.Range("I31").Value2 = 3407.3667
$ws.Range("K31").Value2 = 3407.3667
$ws.Range("M31").Value2 = -3112.3667
$ws.Range("H34").Value2 = 21743362
$ws.Range("I34").Value2 = 3407.3667
$ws.Range("K34").Value2 = 3407.3667
$ws.Range("M34").Value2 = -3205.3667
$ws.Range("H99").Value2 = 5708.0625
$ws.Range("J99").Value2 = 3159.6
$ws.Range("L99").Value2 = 3159.6
$ws.Range("N99").Value2 = -6155.6
$ws.Range("H126").Value2 = 5708.0625
$ws.Range("J126").Value2 = 3159.6
$ws.Range("L126").Value2 = 9478.799999999999
$ws.Range("N126").Value2 = -14418.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 1885.875
$ws.Range("I34").Value2 = 524.25
$ws.Range("J34").Value2 = 3247.5
$ws.Range("K34").Value2 = 1572.75
$ws.Range("L34").Value2 = 9742.5
$ws.Range("M34").Value2 = -1488.75
$ws.Range("N34").Value2 = -9910.5
$ws.Range("H113").Value2 = 4874
$ws.Range("J113").Value2 = 4749
$ws.Range("L113").Value2 = 14247
$ws.Range("N113").Value2 = -18587
$ws.Range("H114").Value2 = 10998.333
$ws.Range("J114").Value2 = 13855.143
$ws.Range("L114").Value2 = 41565.429
$ws.Range("N114").Value2 = -48073.429
$ws.Range("H117").Value2 = 2930.889
$ws.Range("I117").Value2 = 1586.7142
$ws.Range("J117").Value2 = 3786.2727
$ws.Range("K117").Value2 = 4760.142599999999
$ws.Range("L117").Value2 = 11358.8181
$ws.Range("M117").Value2 = -1318.142599999999
$ws.Range("N117").Value2 = -18242.8181
$ws.Range("H129").Value2 = 4081.96
$ws.Range("I129").Value2 = 4303.3335
$ws.Range("J129").Value2 = 3957.4375
$ws.Range("K129").Value2 = 12910.0005
$ws.Range("L129").Value2 = 11872.3125
$ws.Range("M129").Value2 = -7910.000499999998
$ws.Range("N129").Value2 = -21872.3125
$ws.Range("H131").Value2 = 60960.684
$ws.Range("I131").Value2 = 174395
$ws.Range("J131").Value2 = 8606.385
$ws.Range("K131").Value2 = 523185
$ws.Range("L131").Value2 = 25819.155
$ws.Range("M131").Value2 = -518145
$ws.Range("N131").Value2 = -35899.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 12005513
$ws.Range("I126").Value2 = 7149456.5
$ws.Range("J126").Value2 = 18185950
$ws.Range("K126").Value2 = 21448369.5
$ws.Range("L126").Value2 = 54557850
$ws.Range("M126").Value2 = -21445899.5
$ws.Range("N126").Value2 = -54562790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3814.7368
$ws.Range("I7").Value2 = 3696.4707
$ws.Range("J7").Value2 = 4820
$ws.Range("K7").Value2 = 3696.4707
$ws.Range("L7").Value2 = 4820
$ws.Range("M7").Value2 = -3584.4707
$ws.Range("N7").Value2 = -5044
$ws.Range("H40").Value2 = 5244.5
$ws.Range("I40").Value2 = 5061.467
$ws.Range("K40").Value2 = 5061.467
$ws.Range("M40").Value2 = -4925.467
$ws.Range("H46").Value2 = 1250.0358
$ws.Range("I46").Value2 = 899.7619
$ws.Range("J46").Value2 = 2300.8572
$ws.Range("K46").Value2 = 899.7619
$ws.Range("L46").Value2 = 2300.8572
$ws.Range("M46").Value2 = -711.7619
$ws.Range("N46").Value2 = -2676.8572
$ws.Range("H82").Value2 = 5077.222
$ws.Range("I82").Value2 = 3297
$ws.Range("J82").Value2 = 5585.857
$ws.Range("K82").Value2 = 3297
$ws.Range("L82").Value2 = 5585.857
$ws.Range("M82").Value2 = -2936
$ws.Range("N82").Value2 = -6307.857
$ws.Range("H85").Value2 = 5077.222
$ws.Range("I85").Value2 = 3297
$ws.Range("J85").Value2 = 5585.857
$ws.Range("K85").Value2 = 3297
$ws.Range("L85").Value2 = 5585.857
$ws.Range("M85").Value2 = -2049
$ws.Range("N85").Value2 = -8081.857
$ws.Range("H122").Value2 = 5349.875
$ws.Range("I122").Value2 = 5011
$ws.Range("J122").Value2 = 5785.5713
$ws.Range("K122").Value2 = 15033
$ws.Range("L122").Value2 = 17356.7139
$ws.Range("M122").Value2 = -12583
$ws.Range("N122").Value2 = -22256.7139
$ws.Range("H126").Value2 = 3814.7368
$ws.Range("I126").Value2 = 3696.4707
$ws.Range("J126").Value2 = 4820
$ws.Range("K126").Value2 = 11089.4121
$ws.Range("L126").Value2 = 14460
$ws.Range("M126").Value2 = -8619.4121
$ws.Range("N126").Value2 = -19400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value2 = 24999.334
$ws.Range("J59").Value2 = 24999.334
$ws.Range("L59").Value2 = 24999.334
$ws.Range("N59").Value2 = -26475.334
$ws.Range("H136").Value2 = 2556.725
$ws.Range("I136").Value2 = 2478.5
$ws.Range("K136").Value2 = 7435.5
$ws.Range("M136").Value2 = -4885.5
